$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "保險" (Insurance) -- adds company/name/owner header labels
# plus property_category/category/date/legislator_name/legislator_id/
# source_file/index columns to match the common schema used by the
# other property sheets.
# ---------------------------------------------------------------
$wsIns = $wb.Worksheets.Item("保險")

# Fix up the header row (it previously duplicated row 2's data values
# instead of holding field names).
$wsIns.Range("B1").Value = "company"
$wsIns.Range("C1").Value = "name"
$wsIns.Range("D1").Value = "owner"
$wsIns.Range("E1").Value = "property_category"
$wsIns.Range("F1").Value = "category"
$wsIns.Range("G1").Value = "date"
$wsIns.Range("H1").Value = "legislator_name"
$wsIns.Range("I1").Value = "legislator_id"
$wsIns.Range("J1").Value = "source_file"
$wsIns.Range("K1").Value = "index"

# Row 2 (富邦人壽 / 富邦人壽心得意利率變動型年金保險)
$wsIns.Range("E2").Value = "insurance"
$wsIns.Range("F2").Value = "normal"
$wsIns.Range("G2").NumberFormat = "@"
$wsIns.Range("G2").Value = "2013-12-02"
$wsIns.Range("H2").Value = "李俊俋"
$wsIns.Range("I2").Value = 1738
$wsIns.Range("J2").Value = "tmp52b51"
$wsIns.Range("K2").Value = 89

# Row 3 (南山人壽 / 月月金喜利率變動型養老保險)
$wsIns.Range("E3").Value = "insurance"
$wsIns.Range("F3").Value = "normal"
$wsIns.Range("G3").NumberFormat = "@"
$wsIns.Range("G3").Value = "2013-12-02"
$wsIns.Range("H3").Value = "李俊俋"
$wsIns.Range("I3").Value = 1738
$wsIns.Range("J3").Value = "tmp52b51"
$wsIns.Range("K3").Value = 90

# ---------------------------------------------------------------
# Sheet "債務" (Debt) -- adds species/debtor/owner/total header labels
# plus register_date/register_reason/property_category/category/date/
# legislator_name/legislator_id/source_file/index columns.
# ---------------------------------------------------------------
$wsDebt = $wb.Worksheets.Item("債務")

# Fix up the header row.
$wsDebt.Range("B1").Value = "species"
$wsDebt.Range("C1").Value = "debtor"
$wsDebt.Range("D1").Value = "owner"
$wsDebt.Range("E1").Value = "total"
$wsDebt.Range("F1").Value = "register_date"
$wsDebt.Range("G1").Value = "register_reason"
$wsDebt.Range("H1").Value = "property_category"
$wsDebt.Range("I1").Value = "category"
$wsDebt.Range("J1").Value = "date"
$wsDebt.Range("K1").Value = "legislator_name"
$wsDebt.Range("L1").Value = "legislator_id"
$wsDebt.Range("M1").Value = "source_file"
$wsDebt.Range("N1").Value = "index"

# Row 2 (房屋貸款 / 李俊倍 / 第一商業銀行...)
$wsDebt.Range("H2").Value = "debt"
$wsDebt.Range("I2").Value = "normal"
$wsDebt.Range("J2").NumberFormat = "@"
$wsDebt.Range("J2").Value = "2013-12-02"
$wsDebt.Range("K2").Value = "李俊俋"
$wsDebt.Range("L2").Value = 1738
$wsDebt.Range("M2").Value = "tmp52b51"
$wsDebt.Range("N2").Value = 100

# Row 3 (房屋貸款 / 陳佳慧 / 玉山商業銀行...) -- the debt amount had
# been stored as text; correct it to a real number while we are here.
$wsDebt.Range("E3").Value = 3893589
$wsDebt.Range("H3").Value = "debt"
$wsDebt.Range("I3").Value = "normal"
$wsDebt.Range("J3").NumberFormat = "@"
$wsDebt.Range("J3").Value = "2013-12-02"
$wsDebt.Range("K3").Value = "李俊俋"
$wsDebt.Range("L3").Value = 1738
$wsDebt.Range("M3").Value = "tmp52b51"
$wsDebt.Range("N3").Value = 101
